$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 54980
$ws.Range("J126").Value = 54980
$ws.Range("L126").Value = 54980
$ws.Range("N126").Value = -64860
$ws.Range("H132").Value = 478990.06
$ws.Range("I132").Value = 2869.375
$ws.Range("J132").Value = 2002576.2
$ws.Range("K132").Value = 8608.125
$ws.Range("L132").Value = 6007728.6
$ws.Range("M132").Value = -6078.125
$ws.Range("N132").Value = -6012788.6
$ws.Range("H136").Value = 33112.625
$ws.Range("J136").Value = 33112.625
$ws.Range("L136").Value = 33112.625
$ws.Range("N136").Value = -43312.625
$ws.Range("H137").Value = 983.8
$ws.Range("I137").Value = 807.8570999999999
$ws.Range("J137").Value = 1394.3334
$ws.Range("K137").Value = 2423.5713
$ws.Range("L137").Value = 4183.0002
$ws.Range("M137").Value = 126.4287000000004
$ws.Range("N137").Value = -9283.0002
$ws.Range("H138").Value = 2368.2693
$ws.Range("I138").Value = 2175.2307
$ws.Range("J138").Value = 2432.6155
$ws.Range("K138").Value = 6525.6921
$ws.Range("L138").Value = 7297.8465
$ws.Range("M138").Value = -1385.6921
$ws.Range("N138").Value = -17577.8465
$ws.Range("H139").Value = 31093.666
$ws.Range("J139").Value = 31093.666
$ws.Range("L139").Value = 31093.666
$ws.Range("N139").Value = -41373.666
$ws.Range("H141").Value = 2383.9
$ws.Range("I141").Value = 1537.6666
$ws.Range("K141").Value = 4612.9998
$ws.Range("M141").Value = 567.0002000000004

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 28600
$ws.Range("J64").Value = 28600
$ws.Range("L64").Value = 28600
$ws.Range("N64").Value = -29096
$ws.Range("H67").Value = 28600
$ws.Range("J67").Value = 28600
$ws.Range("L67").Value = 28600
$ws.Range("N67").Value = -30316
$ws.Range("H102").Value = 3932.0667
$ws.Range("I102").Value = 3663.3333
$ws.Range("J102").Value = 4335.1665
$ws.Range("K102").Value = 3663.3333
$ws.Range("L102").Value = 4335.1665
$ws.Range("M102").Value = -2041.3333
$ws.Range("N102").Value = -7579.1665
$ws.Range("H132").Value = 2601.75
$ws.Range("I132").Value = 2340.2632
$ws.Range("J132").Value = 3595.4
$ws.Range("K132").Value = 7020.7896
$ws.Range("L132").Value = 10786.2
$ws.Range("M132").Value = -4490.7896
$ws.Range("N132").Value = -15846.2

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 84158.22
$ws.Range("I20").Value = 216529.33
$ws.Range("J20").Value = 17972.666
$ws.Range("K20").Value = 216529.33
$ws.Range("L20").Value = 17972.666
$ws.Range("M20").Value = -216282.33
$ws.Range("N20").Value = -18466.666

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H76").Value = 111115880
$ws.Range("I76").Value = 111115880
$ws.Range("K76").Value = 111115880
$ws.Range("M76").Value = -111115565
$ws.Range("H79").Value = 111115880
$ws.Range("I79").Value = 111115880
$ws.Range("K79").Value = 111115880
$ws.Range("M79").Value = -111114788
$ws.Range("H99").Value = 2519.4285
$ws.Range("I99").Value = 1909
$ws.Range("J99").Value = 3333.3333
$ws.Range("K99").Value = 1909
$ws.Range("L99").Value = 3333.3333
$ws.Range("M99").Value = -411
$ws.Range("N99").Value = -6329.3333
$ws.Range("H109").Value = 22285
$ws.Range("J109").Value = 22285
$ws.Range("L109").Value = 22285
$ws.Range("N109").Value = -24365
$ws.Range("H126").Value = 2519.4285
$ws.Range("I126").Value = 1909
$ws.Range("J126").Value = 3333.3333
$ws.Range("K126").Value = 5727
$ws.Range("L126").Value = 9999.999899999999
$ws.Range("M126").Value = -3257
$ws.Range("N126").Value = -14939.9999
$ws.Range("H138").Value = 42800
$ws.Range("J138").Value = 42800
$ws.Range("L138").Value = 42800
$ws.Range("N138").Value = -53080

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 387.66666
$ws.Range("I92").Value = 390
$ws.Range("J92").Value = 383
$ws.Range("K92").Value = 1170
$ws.Range("L92").Value = 1149
$ws.Range("M92").Value = 78
$ws.Range("N92").Value = -3645
$ws.Range("H137").Value = 46800.793
$ws.Range("I137").Value = 1311.8182
$ws.Range("J137").Value = 85291.46000000001
$ws.Range("K137").Value = 3935.4546
$ws.Range("L137").Value = 255874.38
$ws.Range("M137").Value = 1164.5454
$ws.Range("N137").Value = -266074.38

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3644.7368
$ws.Range("I80").Value = 3573.5293
$ws.Range("J80").Value = 4250
$ws.Range("K80").Value = 3573.5293
$ws.Range("L80").Value = 4250
$ws.Range("M80").Value = -2575.5293
$ws.Range("N80").Value = -6246
$ws.Range("H83").Value = 3644.7368
$ws.Range("I83").Value = 3573.5293
$ws.Range("J83").Value = 4250
$ws.Range("K83").Value = 17867.6465
$ws.Range("L83").Value = 21250
$ws.Range("M83").Value = -12875.6465
$ws.Range("N83").Value = -31234
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H136").Value = 23000.25
$ws.Range("J136").Value = 23000.25
$ws.Range("L136").Value = 69000.75
$ws.Range("N136").Value = -74100.75
$ws.Range("H139").Value = 23303
$ws.Range("J139").Value = 23303
$ws.Range("L139").Value = 23303
$ws.Range("N139").Value = -33583

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2024.75
$ws.Range("I7").Value = 2050
$ws.Range("J7").Value = 1999.5
$ws.Range("K7").Value = 2050
$ws.Range("L7").Value = 1999.5
$ws.Range("M7").Value = -1938
$ws.Range("N7").Value = -2223.5
$ws.Range("H40").Value = 2050
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
$ws.Range("H68").Value = 2497.6365
$ws.Range("I68").Value = 2182
$ws.Range("J68").Value = 3050
$ws.Range("K68").Value = 2182
$ws.Range("L68").Value = 3050
$ws.Range("M68").Value = -1433
$ws.Range("N68").Value = -4548
$ws.Range("H71").Value = 2497.6365
$ws.Range("I71").Value = 2182
$ws.Range("J71").Value = 3050
$ws.Range("K71").Value = 10910
$ws.Range("L71").Value = 15250
$ws.Range("M71").Value = -7166
$ws.Range("N71").Value = -22738
$ws.Range("H126").Value = 2024.75
$ws.Range("I126").Value = 2050
$ws.Range("J126").Value = 1999.5
$ws.Range("K126").Value = 6150
$ws.Range("L126").Value = 5998.5
$ws.Range("M126").Value = -3680
$ws.Range("N126").Value = -10938.5
$ws.Range("H134").Value = 56000
$ws.Range("J134").Value = 56000
$ws.Range("L134").Value = 56000
$ws.Range("N134").Value = -66140
$ws.Range("H135").Value = 55771.43
$ws.Range("J135").Value = 55771.43
$ws.Range("L135").Value = 55771.43
$ws.Range("N135").Value = -65911.42999999999
$ws.Range("H138").Value = 30985.8
$ws.Range("J138").Value = 30985.8
$ws.Range("L138").Value = 30985.8
$ws.Range("N138").Value = -41265.8

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3416.3333
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 3624.5
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 10873.5
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -15933.5
$ws.Range("H138").Value = 54000
$ws.Range("J138").Value = 54000
$ws.Range("L138").Value = 54000
$ws.Range("N138").Value = -64280
